# Applies cryptocurrency price/volume updates per commit "Updated cryptos list ... with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.183.63"
$ws.Range("E2").Value = "  +0.75%  "

$ws.Range("D3").Value = "3.503.20"
$ws.Range("E3").Value = "  +0.03%  "

$ws.Range("E4").Value = "  +0.00%  "

$ws.Range("D5").Value = "'599.34"
$ws.Range("E5").Value = "  +1.00%  "

$ws.Range("D6").Value = "'173.94"
$ws.Range("E6").Value = "  +2.75%  "

$ws.Range("E7").Value = "  -0.02%  "

$ws.Range("D8").Value = "'0.586"
$ws.Range("E8").Value = "  -1.05%  "

$ws.Range("E9").Value = "  -0.61%  "

$ws.Range("D10").Value = "'7.18"
$ws.Range("E10").Value = "  -2.10%  "

$ws.Range("D11").Value = "'0.431"
$ws.Range("E11").Value = "  -0.48%  "

$ws.Range("D12").Value = "4.110.79"
$ws.Range("E12").Value = "  +0.11%  "

$ws.Range("D13").Value = "'30.62"
$ws.Range("E13").Value = "  +8.21%  "

$ws.Range("E14").Value = "  +0.28%  "

$ws.Range("D15").Value = "67.148.32"
$ws.Range("E15").Value = "  +0.70%  "

$ws.Range("D16").Value = "'0.0000179"
$ws.Range("E16").Value = "  -1.67%  "

$ws.Range("D17").Value = "3.494.35"
$ws.Range("E17").Value = "  +0.06%  "

$ws.Range("D18").Value = "'6.30"
$ws.Range("E18").Value = "  -0.44%  "

$ws.Range("D19").Value = "'14.61"
$ws.Range("E19").Value = "  +3.96%  "

$ws.Range("D20").Value = "'393.94"
$ws.Range("E20").Value = "  -0.78%  "

$ws.Range("D21").Value = "'7.98"
$ws.Range("E21").Value = "  +0.27%  "

$ws.Range("D22").Value = "'73.36"
$ws.Range("E22").Value = "  -0.05%  "

$ws.Range("D23").Value = "'0.998"
$ws.Range("E23").Value = "  -0.12%  "

$ws.Range("D24").Value = "'0.537"
$ws.Range("E24").Value = "  +0.46%  "

$ws.Range("E25").Value = "  -0.82%  "

$ws.Range("E26").Value = "  -0.38%  "

$ws.Range("D27").Value = "'10.13"
$ws.Range("E27").Value = "  -0.76%  "

$ws.Range("E28").Value = "  -0.12%  "

$ws.Range("D29").Value = "'0.995"
$ws.Range("E29").Value = "  -0.41%  "

$ws.Range("D30").Value = "'6.16"
$ws.Range("E30").Value = "  -2.56%  "

$ws.Range("D31").Value = "'1.42"
$ws.Range("E31").Value = "  -2.61%  "

$ws.Range("D32").Value = "'2.06"
$ws.Range("E32").Value = "  -0.23%  "

$ws.Range("D33").Value = "'23.64"
$ws.Range("E33").Value = "  -0.69%  "

$ws.Range("D34").Value = "'7.39"
$ws.Range("E34").Value = "  -0.27%  "

$ws.Range("D35").Value = "'1.63"
$ws.Range("E35").Value = "  +1.35%  "

$ws.Range("D36").Value = "'164.04"
$ws.Range("E36").Value = "  +1.07%  "

$ws.Range("D37").Value = "'0.878"
$ws.Range("E37").Value = "  -2.77%  "

$ws.Range("E38").Value = "  -0.01%  "

$ws.Range("D39").Value = "'7.02"
$ws.Range("E39").Value = "  +3.27%  "

$ws.Range("D40").Value = "'4.68"
$ws.Range("E40").Value = "  -0.09%  "

$ws.Range("D41").Value = "'27.46"
$ws.Range("E41").Value = "  +0.61%  "

$ws.Range("D42").Value = "'0.0731"
$ws.Range("E42").Value = "  -1.90%  "

$ws.Range("D43").Value = "'26.10"
$ws.Range("E43").Value = "  -1.72%  "

$ws.Range("D44").Value = "2.802.11"
$ws.Range("E44").Value = "  -0.32%  "

$ws.Range("D45").Value = "'42.53"
$ws.Range("E45").Value = "  -1.00%  "

$ws.Range("D46").Value = "'2.53"
$ws.Range("E46").Value = "  -1.41%  "

$ws.Range("D47").Value = "'0.0303"
$ws.Range("E47").Value = "  -3.35%  "

$ws.Range("D48").Value = "'341.27"
$ws.Range("E48").Value = "  -0.54%  "

$ws.Range("D49").Value = "'1.08"
$ws.Range("E49").Value = "  -1.53%  "

$ws.Range("D50").Value = "'33.65"
$ws.Range("E50").Value = "  -1.47%  "

$ws.Range("B51").Value = "Cosmos"
$ws.Range("C51").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D51").Value = "'6.44"
$ws.Range("E51").Value = "  -1.06%  "
